$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 933.6
$ws.Range("I12").Value = 334
$ws.Range("K12").Value = 334
$ws.Range("M12").Value = -164
$ws.Range("H33").Value = 50000292
$ws.Range("I33").Value = 83333450
$ws.Range("J33").Value = 559
$ws.Range("K33").Value = 83333450
$ws.Range("L33").Value = 559
$ws.Range("M33").Value = -83333221
$ws.Range("N33").Value = -1017
$ws.Range("H87").Value = 61000
$ws.Range("J87").Value = 61000
$ws.Range("L87").Value = 61000
$ws.Range("N87").Value = -63496
$ws.Range("H90").Value = 61000
$ws.Range("J90").Value = 61000
$ws.Range("L90").Value = 183000
$ws.Range("N90").Value = -195480
$ws.Range("H98").Value = 1287.4615
$ws.Range("I98").Value = 1019.8333
$ws.Range("K98").Value = 1019.8333
$ws.Range("M98").Value = 478.1667
$ws.Range("H106").Value = 2641.7144
$ws.Range("I106").Value = 2590.3333
$ws.Range("K106").Value = 2590.3333
$ws.Range("M106").Value = -1959.3333
$ws.Range("H122").Value = 1287.4615
$ws.Range("I122").Value = 1019.8333
$ws.Range("K122").Value = 3059.4999
$ws.Range("M122").Value = -609.4998999999998
$ws.Range("H131").Value = 127287
$ws.Range("I131").Value = 202519.2
$ws.Range("K131").Value = 607557.6000000001
$ws.Range("M131").Value = -602517.6000000001
$ws.Range("H132").Value = 3295.8
$ws.Range("I132").Value = 882.1212
$ws.Range("K132").Value = 2646.3636
$ws.Range("M132").Value = -116.3636000000001
$ws.Range("H135").Value = 563.9375
$ws.Range("I135").Value = 504.92307
$ws.Range("K135").Value = 4544.30763
$ws.Range("M135").Value = -2009.30763
$ws.Range("H137").Value = 2914.889
$ws.Range("I137").Value = 2889.1667
$ws.Range("K137").Value = 8667.500100000001
$ws.Range("M137").Value = -6117.500100000001

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17032.453
$ws.Range("I32").Value = 3759.2
$ws.Range("K32").Value = 3759.2
$ws.Range("M32").Value = -3472.2
$ws.Range("H45").Value = 3772.6667
$ws.Range("I45").Value = 2673.2222
$ws.Range("J45").Value = 4432.3335
$ws.Range("K45").Value = 2673.2222
$ws.Range("L45").Value = 4432.3335
$ws.Range("M45").Value = -2296.2222
$ws.Range("N45").Value = -5186.3335
$ws.Range("H61").Value = 11372834
$ws.Range("I61").Value = 17247974
$ws.Range("J61").Value = 14227.4
$ws.Range("K61").Value = 17247974
$ws.Range("L61").Value = 14227.4
$ws.Range("M61").Value = -17247762
$ws.Range("N61").Value = -14651.4
$ws.Range("H74").Value = 4422.294
$ws.Range("I74").Value = 3055.1538
$ws.Range("J74").Value = 8865.5
$ws.Range("K74").Value = 3055.1538
$ws.Range("L74").Value = 8865.5
$ws.Range("M74").Value = -2181.1538
$ws.Range("N74").Value = -10613.5
$ws.Range("H77").Value = 4422.294
$ws.Range("I77").Value = 3055.1538
$ws.Range("J77").Value = 8865.5
$ws.Range("K77").Value = 15275.769
$ws.Range("L77").Value = 44327.5
$ws.Range("M77").Value = -10907.769
$ws.Range("N77").Value = -53063.5
$ws.Range("H110").Value = 4762.4116
$ws.Range("I110").Value = 3812.111
$ws.Range("J110").Value = 8427.857
$ws.Range("K110").Value = 3812.111
$ws.Range("L110").Value = 8427.857
$ws.Range("M110").Value = -1767.111
$ws.Range("N110").Value = -12517.857
$ws.Range("H122").Value = 1559.375
$ws.Range("I122").Value = 1353.5714
$ws.Range("K122").Value = 4060.7142
$ws.Range("M122").Value = -1610.7142
$ws.Range("H132").Value = 2006.38
$ws.Range("I132").Value = 2149.2444
$ws.Range("K132").Value = 6447.733200000001
$ws.Range("M132").Value = -3917.733200000001
$ws.Range("H136").Value = 11372834
$ws.Range("I136").Value = 17247974
$ws.Range("J136").Value = 14227.4
$ws.Range("K136").Value = 51743922
$ws.Range("L136").Value = 42682.2
$ws.Range("M136").Value = -51741372
$ws.Range("N136").Value = -47782.2

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 442.375
$ws.Range("J11").Value = 767
$ws.Range("L11").Value = 767
$ws.Range("N11").Value = -1047
$ws.Range("H20").Value = 4996.3335
$ws.Range("I20").Value = 4676.727
$ws.Range("K20").Value = 4676.727
$ws.Range("M20").Value = -4429.727
$ws.Range("H105").Value = 1615.5
$ws.Range("I105").Value = 1473.875
$ws.Range("K105").Value = 1473.875
$ws.Range("M105").Value = 273.125
$ws.Range("H107").Value = 5725.727
$ws.Range("I107").Value = 5426.857
$ws.Range("K107").Value = 5426.857
$ws.Range("M107").Value = -3506.857
$ws.Range("H134").Value = 2402.0142
$ws.Range("I134").Value = 2518.0952
$ws.Range("K134").Value = 7554.285600000001
$ws.Range("M134").Value = -5019.285600000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3427.6875
$ws.Range("I16").Value = 1955.5
$ws.Range("K16").Value = 1955.5
$ws.Range("M16").Value = -1668.5
$ws.Range("H31").Value = 6043.4707
$ws.Range("I31").Value = 3333.111
$ws.Range("K31").Value = 3333.111
$ws.Range("M31").Value = -3038.111
$ws.Range("H34").Value = 6043.4707
$ws.Range("I34").Value = 3333.111
$ws.Range("K34").Value = 3333.111
$ws.Range("M34").Value = -3131.111
$ws.Range("H97").Value = 38000
$ws.Range("J97").Value = 38000
$ws.Range("L97").Value = 38000
$ws.Range("N97").Value = -39982
$ws.Range("H99").Value = 3659.75
$ws.Range("I99").Value = 3585.5
$ws.Range("K99").Value = 3585.5
$ws.Range("M99").Value = -2087.5
$ws.Range("H113").Value = 3427.6875
$ws.Range("I113").Value = 1955.5
$ws.Range("K113").Value = 1955.5
$ws.Range("M113").Value = 214.5
$ws.Range("H126").Value = 3659.75
$ws.Range("I126").Value = 3585.5
$ws.Range("K126").Value = 10756.5
$ws.Range("M126").Value = -8286.5
$ws.Range("H132").Value = 3574.476
$ws.Range("I132").Value = 3725.7778
$ws.Range("K132").Value = 11177.3334
$ws.Range("M132").Value = -8647.3334
$ws.Range("H133").Value = 61808.668
$ws.Range("J133").Value = 61808.668
$ws.Range("L133").Value = 61808.668
$ws.Range("N133").Value = -66868.66800000001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1500
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 1500
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 4500
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -5088
$ws.Range("H98").Value = 200.09091
$ws.Range("J98").Value = 200.09091
$ws.Range("L98").Value = 600.27273
$ws.Range("N98").Value = -3596.27273
$ws.Range("H109").Value = 16002.667
$ws.Range("I109").Value = 18803.2
$ws.Range("J109").Value = 2000
$ws.Range("K109").Value = 56409.60000000001
$ws.Range("L109").Value = 6000
$ws.Range("M109").Value = -55369.60000000001
$ws.Range("N109").Value = -8080
$ws.Range("H122").Value = 914
$ws.Range("J122").Value = 1099.3
$ws.Range("L122").Value = 9893.699999999999
$ws.Range("N122").Value = -14793.7

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14433.066
$ws.Range("I70").Value = 10812.25
$ws.Range("K70").Value = 10812.25
$ws.Range("M70").Value = -10542.25
$ws.Range("H73").Value = 14433.066
$ws.Range("I73").Value = 10812.25
$ws.Range("K73").Value = 10812.25
$ws.Range("M73").Value = -9876.25
$ws.Range("H97").Value = 2346.9473
$ws.Range("I97").Value = 884.8
$ws.Range("J97").Value = 5158.769
$ws.Range("K97").Value = 884.8
$ws.Range("L97").Value = 5158.769
$ws.Range("M97").Value = -388.8
$ws.Range("N97").Value = -6150.769
$ws.Range("H132").Value = 1233.8572
$ws.Range("I132").Value = 1145.4916
$ws.Range("K132").Value = 3436.4748
$ws.Range("M132").Value = -906.4748

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2606.158
$ws.Range("I122").Value = 2412.7646
$ws.Range("K122").Value = 7238.293799999999
$ws.Range("M122").Value = -4788.293799999999
$ws.Range("H132").Value = 10033.435
$ws.Range("I132").Value = 10282
$ws.Range("J132").Value = 9138.6
$ws.Range("K132").Value = 30846
$ws.Range("L132").Value = 27415.8
$ws.Range("M132").Value = -28316
$ws.Range("N132").Value = -32475.8
$ws.Range("H136").Value = 4845.5713
$ws.Range("I136").Value = 4950.963
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 14852.889
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -12302.889
$ws.Range("N136").Value = -11100

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1632.5555
$ws.Range("J107").Value = 1949.3334
$ws.Range("L107").Value = 5848.0002
$ws.Range("N107").Value = -9688.0002
$ws.Range("H122").Value = 2652.2964
$ws.Range("J122").Value = 6285.4287
$ws.Range("L122").Value = 18856.2861
$ws.Range("N122").Value = -23756.2861
$ws.Range("H126").Value = 2863.5757
$ws.Range("J126").Value = 3954.25
$ws.Range("L126").Value = 11862.75
$ws.Range("N126").Value = -16802.75
$ws.Range("H132").Value = 1906.7234
$ws.Range("I132").Value = 2018.6136
$ws.Range("J132").Value = 265.66666
$ws.Range("K132").Value = 6055.8408
$ws.Range("L132").Value = 796.9999799999999
$ws.Range("M132").Value = -3525.8408
$ws.Range("N132").Value = -5856.99998
